$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 169.915657
$ws.Cells.Item(2, 8).Value = 509.746971
$ws.Cells.Item(2, 9).Value = 0.4441184931734509
$ws.Cells.Item(2, 10).Value = 0.4441184931734509
$ws.Cells.Item(2, 13).Value = 0.817551
$ws.Cells.Item(2, 14).Value = 2.452653
$ws.Cells.Item(2, 15).Value = 0.2022370099583455
$ws.Cells.Item(2, 16).Value = 0.2022370099583455
$ws.Cells.Item(2, 17).Value = 138.914715296007
$ws.Cells.Item(2, 18).Value = 1250.232437664063
$ws.Cells.Item(2, 19).Value = 0.08981719612660459
$ws.Cells.Item(2, 20).Value = 0.08981719612660459

# Row 3
$ws.Cells.Item(3, 7).Value = 169.915657
$ws.Cells.Item(3, 8).Value = 509.746971
$ws.Cells.Item(3, 9).Value = 0.4441184931734509
$ws.Cells.Item(3, 10).Value = 0.4441184931734509
$ws.Cells.Item(3, 15).Value = 0.2955761218382804
$ws.Cells.Item(3, 16).Value = 0.2955761218382804
$ws.Cells.Item(3, 17).Value = 203.028480404846
$ws.Cells.Item(3, 18).Value = 1827.256323643614
$ws.Cells.Item(3, 19).Value = 0.1312708218488694
$ws.Cells.Item(3, 20).Value = 0.1312708218488694

# Row 4
$ws.Cells.Item(4, 7).Value = 169.915657
$ws.Cells.Item(4, 8).Value = 509.746971
$ws.Cells.Item(4, 9).Value = 0.4441184931734509
$ws.Cells.Item(4, 10).Value = 0.4441184931734509
$ws.Cells.Item(4, 13).Value = 1.814372333333333
$ws.Cells.Item(4, 14).Value = 5.443117
$ws.Cells.Item(4, 15).Value = 0.4488199948926487
$ws.Cells.Item(4, 16).Value = 0.4488199948926487
$ws.Cells.Item(4, 17).Value = 308.2902670609564
$ws.Cells.Item(4, 18).Value = 2774.612403548607
$ws.Cells.Item(4, 19).Value = 0.1993292598378391
$ws.Cells.Item(4, 20).Value = 0.1993292598378391

# Row 5
$ws.Cells.Item(5, 7).Value = 169.915657
$ws.Cells.Item(5, 8).Value = 509.746971
$ws.Cells.Item(5, 9).Value = 0.4441184931734509
$ws.Cells.Item(5, 10).Value = 0.4441184931734509
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.2157376666666667
$ws.Cells.Item(5, 14).Value = 0.647213
$ws.Cells.Item(5, 15).Value = 0.05336687331072544
$ws.Cells.Item(5, 16).Value = 0.05336687331072543
$ws.Cells.Item(5, 17).Value = 36.65720737131367
$ws.Cells.Item(5, 18).Value = 329.914866341823
$ws.Cells.Item(5, 19).Value = 0.02370121536013783
$ws.Cells.Item(5, 20).Value = 0.02370121536013783

# Row 6
$ws.Cells.Item(6, 9).Value = 0.1787346690539575
$ws.Cells.Item(6, 10).Value = 0.1787346690539575
$ws.Cells.Item(6, 13).Value = 0.817551
$ws.Cells.Item(6, 14).Value = 2.452653
$ws.Cells.Item(6, 15).Value = 0.2022370099583455
$ws.Cells.Item(6, 16).Value = 0.2022370099583455
$ws.Cells.Item(6, 17).Value = 55.905971146893
$ws.Cells.Item(6, 18).Value = 503.153740322037
$ws.Cells.Item(6, 19).Value = 0.0361467650453668
$ws.Cells.Item(6, 20).Value = 0.03614676504536679

# Row 7
$ws.Cells.Item(7, 9).Value = 0.1787346690539575
$ws.Cells.Item(7, 10).Value = 0.1787346690539575
$ws.Cells.Item(7, 15).Value = 0.2955761218382804
$ws.Cells.Item(7, 16).Value = 0.2955761218382804
$ws.Cells.Item(7, 19).Value = 0.05282970031701727
$ws.Cells.Item(7, 20).Value = 0.05282970031701727

# Row 8
$ws.Cells.Item(8, 9).Value = 0.1787346690539575
$ws.Cells.Item(8, 10).Value = 0.1787346690539575
$ws.Cells.Item(8, 13).Value = 1.814372333333333
$ws.Cells.Item(8, 14).Value = 5.443117
$ws.Cells.Item(8, 15).Value = 0.4488199948926487
$ws.Cells.Item(8, 16).Value = 0.4488199948926487
$ws.Cells.Item(8, 17).Value = 124.070849790477
$ws.Cells.Item(8, 18).Value = 1116.637648114293
$ws.Cells.Item(8, 19).Value = 0.08021969325193648
$ws.Cells.Item(8, 20).Value = 0.08021969325193648

# Row 9
$ws.Cells.Item(9, 9).Value = 0.1787346690539575
$ws.Cells.Item(9, 10).Value = 0.1787346690539575
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.2157376666666667
$ws.Cells.Item(9, 14).Value = 0.647213
$ws.Cells.Item(9, 15).Value = 0.05336687331072544
$ws.Cells.Item(9, 16).Value = 0.05336687331072543
$ws.Cells.Item(9, 17).Value = 14.752625546253
$ws.Cells.Item(9, 18).Value = 132.773629916277
$ws.Cells.Item(9, 19).Value = 0.00953851043963699
$ws.Cells.Item(9, 20).Value = 0.009538510439636987

# Row 10
$ws.Cells.Item(10, 7).Value = 53.27463399999999
$ws.Cells.Item(10, 8).Value = 159.823902
$ws.Cells.Item(10, 9).Value = 0.1392470275793777
$ws.Cells.Item(10, 10).Value = 0.1392470275793778
$ws.Cells.Item(10, 13).Value = 0.817551
$ws.Cells.Item(10, 14).Value = 2.452653
$ws.Cells.Item(10, 15).Value = 0.2022370099583455
$ws.Cells.Item(10, 16).Value = 0.2022370099583455
$ws.Cells.Item(10, 17).Value = 43.554730301334
$ws.Cells.Item(10, 18).Value = 391.992572712006
$ws.Cells.Item(10, 19).Value = 0.02816090250324063
$ws.Cells.Item(10, 20).Value = 0.02816090250324063

# Row 11
$ws.Cells.Item(11, 7).Value = 53.27463399999999
$ws.Cells.Item(11, 8).Value = 159.823902
$ws.Cells.Item(11, 9).Value = 0.1392470275793777
$ws.Cells.Item(11, 10).Value = 0.1392470275793778
$ws.Cells.Item(11, 15).Value = 0.2955761218382804
$ws.Cells.Item(11, 16).Value = 0.2955761218382804
$ws.Cells.Item(11, 17).Value = 63.65668812465199
$ws.Cells.Item(11, 18).Value = 572.910193121868
$ws.Cells.Item(11, 19).Value = 0.04115809638942055
$ws.Cells.Item(11, 20).Value = 0.04115809638942055

# Row 12
$ws.Cells.Item(12, 7).Value = 53.27463399999999
$ws.Cells.Item(12, 8).Value = 159.823902
$ws.Cells.Item(12, 9).Value = 0.1392470275793777
$ws.Cells.Item(12, 10).Value = 0.1392470275793778
$ws.Cells.Item(12, 13).Value = 1.814372333333333
$ws.Cells.Item(12, 14).Value = 5.443117
$ws.Cells.Item(12, 15).Value = 0.4488199948926487
$ws.Cells.Item(12, 16).Value = 0.4488199948926487
$ws.Cells.Item(12, 17).Value = 96.66002199805932
$ws.Cells.Item(12, 18).Value = 869.9401979825338
$ws.Cells.Item(12, 19).Value = 0.06249685020699283
$ws.Cells.Item(12, 20).Value = 0.06249685020699285

# Row 13
$ws.Cells.Item(13, 7).Value = 53.27463399999999
$ws.Cells.Item(13, 8).Value = 159.823902
$ws.Cells.Item(13, 9).Value = 0.1392470275793777
$ws.Cells.Item(13, 10).Value = 0.1392470275793778
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.2157376666666667
$ws.Cells.Item(13, 14).Value = 0.647213
$ws.Cells.Item(13, 15).Value = 0.05336687331072544
$ws.Cells.Item(13, 16).Value = 0.05336687331072543
$ws.Cells.Item(13, 17).Value = 11.49334523168067
$ws.Cells.Item(13, 18).Value = 103.440107085126
$ws.Cells.Item(13, 19).Value = 0.007431178479723743
$ws.Cells.Item(13, 20).Value = 0.007431178479723743

# Row 14
$ws.Cells.Item(14, 7).Value = 91.01828266666666
$ws.Cells.Item(14, 8).Value = 273.054848
$ws.Cells.Item(14, 9).Value = 0.2378998101932138
$ws.Cells.Item(14, 10).Value = 0.2378998101932138
$ws.Cells.Item(14, 13).Value = 0.817551
$ws.Cells.Item(14, 14).Value = 2.452653
$ws.Cells.Item(14, 15).Value = 0.2022370099583455
$ws.Cells.Item(14, 16).Value = 0.2022370099583455
$ws.Cells.Item(14, 17).Value = 74.412088012416
$ws.Cells.Item(14, 18).Value = 669.7087921117441
$ws.Cells.Item(14, 19).Value = 0.04811214628313348
$ws.Cells.Item(14, 20).Value = 0.04811214628313348

# Row 15
$ws.Cells.Item(15, 7).Value = 91.01828266666666
$ws.Cells.Item(15, 8).Value = 273.054848
$ws.Cells.Item(15, 9).Value = 0.2378998101932138
$ws.Cells.Item(15, 10).Value = 0.2378998101932138
$ws.Cells.Item(15, 15).Value = 0.2955761218382804
$ws.Cells.Item(15, 16).Value = 0.2955761218382804
$ws.Cells.Item(15, 17).Value = 108.7557435561813
$ws.Cells.Item(15, 18).Value = 978.801692005632
$ws.Cells.Item(15, 19).Value = 0.07031750328297315
$ws.Cells.Item(15, 20).Value = 0.07031750328297313

# Row 16
$ws.Cells.Item(16, 7).Value = 91.01828266666666
$ws.Cells.Item(16, 8).Value = 273.054848
$ws.Cells.Item(16, 9).Value = 0.2378998101932138
$ws.Cells.Item(16, 10).Value = 0.2378998101932138
$ws.Cells.Item(16, 13).Value = 1.814372333333333
$ws.Cells.Item(16, 14).Value = 5.443117
$ws.Cells.Item(16, 15).Value = 0.4488199948926487
$ws.Cells.Item(16, 16).Value = 0.4488199948926487
$ws.Cells.Item(16, 17).Value = 165.1410538979129
$ws.Cells.Item(16, 18).Value = 1486.269485081216
$ws.Cells.Item(16, 19).Value = 0.1067741915958803
$ws.Cells.Item(16, 20).Value = 0.1067741915958803

# Row 17
$ws.Cells.Item(17, 7).Value = 91.01828266666666
$ws.Cells.Item(17, 8).Value = 273.054848
$ws.Cells.Item(17, 9).Value = 0.2378998101932138
$ws.Cells.Item(17, 10).Value = 0.2378998101932138
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.2157376666666667
$ws.Cells.Item(17, 14).Value = 0.647213
$ws.Cells.Item(17, 15).Value = 0.05336687331072544
$ws.Cells.Item(17, 16).Value = 0.05336687331072543
$ws.Cells.Item(17, 17).Value = 19.63607192651378
$ws.Cells.Item(17, 18).Value = 176.724647338624
$ws.Cells.Item(17, 19).Value = 0.01269596903122687
$ws.Cells.Item(17, 20).Value = 0.01269596903122687
